$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 42, pushing the existing rows 42-52 down to 44-54.
$ws.Rows.Item(42).Insert()
$ws.Rows.Item(42).Insert()

# --- New row 42: "Dina" / "Especial" ---
$ws.Range("A42").Value = 9
$ws.Range("B42").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C42").Value = "Metropolitana"
$ws.Range("D42").Value = 44551
$ws.Range("D42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E42").Value = 13
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100103
$ws.Range("H42").Value = "Frutos de hueso (carozo)"
$ws.Range("I42").Value = 100103003
$ws.Range("J42").Value = "Damasco"
$ws.Range("K42").Value = "Dina"
$ws.Range("L42").Value = "Especial"
$ws.Range("M42").Value = 220
$ws.Range("N42").Value = 18000
$ws.Range("O42").Value = 18000
$ws.Range("P42").Value = 18000
$ws.Range("Q42").Value = "$/caja 18 kilos granel"
$ws.Range("R42").Value = "Provincia de Los Andes"
$ws.Range("S42").Value = 1000
$ws.Range("T42").Value = 18

# --- New row 43: "Dina" / "Primera" ---
$ws.Range("A43").Value = 9
$ws.Range("B43").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C43").Value = "Metropolitana"
$ws.Range("D43").Value = 44551
$ws.Range("D43").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E43").Value = 13
$ws.Range("F43").Value = "Fruta"
$ws.Range("G43").Value = 100103
$ws.Range("H43").Value = "Frutos de hueso (carozo)"
$ws.Range("I43").Value = 100103003
$ws.Range("J43").Value = "Damasco"
$ws.Range("K43").Value = "Dina"
$ws.Range("L43").Value = "Primera"
$ws.Range("M43").Value = 350
$ws.Range("N43").Value = 14400
$ws.Range("O43").Value = 14400
$ws.Range("P43").Value = 14400
$ws.Range("Q43").Value = "$/caja 18 kilos granel"
$ws.Range("R43").Value = "Provincia de Los Andes"
$ws.Range("S43").Value = 800
$ws.Range("T43").Value = 18
